$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "Aircraft ActiveTrack available at max speed . When exceeding nnn, Obstacle Avoidance is not available ."
$ws.Range("C2").Value = "When exceeding nnn, Obstacle Avoidance is not available"
$ws.Range("D2").Value = "7-14"
$ws.Range("E2").Value = "Missing"

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "Aircraft ActiveTrack available at max speed . When exceeding nnn, Obstacle Avoidance is not available ."
$ws.Range("C3").Value = "Obstacle Avoidance is not available"
$ws.Range("D3").Value = "10-14"
$ws.Range("E3").Value = "'False"

$ws.Range("A4").Value = 19
$ws.Range("B4").Value = "Aircraft is tilted , please keep the aircraft stationary and level before flight ."
$ws.Range("C4").Value = "Aircraft is tilted"
$ws.Range("D4").Value = "0-2"
$ws.Range("E4").Value = "Missing"

$ws.Range("A5").Value = 40
$ws.Range("B5").Value = "Camera error . AI Spot-Check failed . Restart camera ."
$ws.Range("C5").Value = "AI Spot-Check failed"
$ws.Range("D5").Value = "3-5"
$ws.Range("E5").Value = "Missing"

$ws.Range("A6").Value = 42
$ws.Range("B6").Value = "Camera sensor error . Hardware malfunction : Contact DJI Support to arrange for repairs ."
$ws.Range("C6").Value = "Hardware malfunction"
$ws.Range("D6").Value = "4-5"
$ws.Range("E6").Value = "Missing"

$ws.Range("A7").Value = 42
$ws.Range("B7").Value = "Camera sensor error . Hardware malfunction : Contact DJI Support to arrange for repairs ."
$ws.Range("C7").Value = "Contact DJI Support to arrange for repairs"
$ws.Range("D7").Value = "7-13"
$ws.Range("E7").Value = "Missing"

$ws.Range("A8").Value = 48
$ws.Range("B8").Value = "Check whether propellers are installed correctly . If the propellers are installed correctly and the aircraft still cannot takeoff, a motor error may exist . Contact DJI Support for assistance ."
$ws.Range("C8").Value = "If the propellers are installed correctly and the aircraft still cannot takeoff, a motor error may exist"
$ws.Range("D8").Value = "7-23"
$ws.Range("E8").Value = "Missing"

$ws.Range("A9").Value = 48
$ws.Range("B9").Value = "Check whether propellers are installed correctly . If the propellers are installed correctly and the aircraft still cannot takeoff, a motor error may exist . Contact DJI Support for assistance ."
$ws.Range("C9").Value = "If the propellers are installed correctly and the aircraft still cannot takeoff,"
$ws.Range("D9").Value = "7-18"
$ws.Range("E9").Value = "Missing"

$ws.Range("A10").Value = 50
$ws.Range("B10").Value = "Compass abnormal . Solution: 1. Ensure there are no magnets or metal objects near the aircraft . The ground or walls may contain metal . Move away from sources of interference before attempting flight . 2. Calibrate Compass Before Takeoff ."
$ws.Range("C10").Value = "2. Calibrate Compass Before Takeoff"
$ws.Range("D10").Value = "35-39"
$ws.Range("E10").Value = "Missing"

$ws.Range("A11").Value = 50
$ws.Range("B11").Value = "Compass abnormal . Solution: 1. Ensure there are no magnets or metal objects near the aircraft . The ground or walls may contain metal . Move away from sources of interference before attempting flight . 2. Calibrate Compass Before Takeoff ."
$ws.Range("C11").Value = "Calibrate Compass Before Takeoff"
$ws.Range("D11").Value = "36-39"
$ws.Range("E11").Value = "'False"

$ws.Range("A12").Value = 66
$ws.Range("B12").Value = "Downlink data connection lost for nnn seconds ."
$ws.Range("C12").Value = "Downlink data connection lost for nnn seconds"
$ws.Range("D12").Value = "0-6"
$ws.Range("E12").Value = "Missing"

$ws.Range("A13").Value = 66
$ws.Range("B13").Value = "Downlink data connection lost for nnn seconds ."
$ws.Range("C13").Value = "Downlink data connection lost for nnn"
$ws.Range("D13").Value = "0-5"
$ws.Range("E13").Value = "'False"

$ws.Range("A14").Value = 77
$ws.Range("B14").Value = "Exiting GPS mode : Unknown Error ."
$ws.Range("C14").Value = "Unknown Error"
$ws.Range("D14").Value = "4-5"
$ws.Range("E14").Value = "Missing"

$ws.Range("A15").Value = 81
$ws.Range("B15").Value = "Extra payload detected . Return aircraft to an area nearby the home point promptly and fly in a wind-free environment to ensure flight safety ."
$ws.Range("C15").Value = "Return aircraft to an area nearby the home point promptly and fly in a wind-free environment to ensure flight safety"
$ws.Range("D15").Value = "4-23"
$ws.Range("E15").Value = "Missing"

$ws.Range("A16").Value = 85
$ws.Range("B16").Value = "Flight altitude exceeds nnn . Aircraft may be in violation of local laws and regulations . Check and make sure you have obtained proper authorization to fly in this airspace ."
$ws.Range("C16").Value = "Aircraft may be in violation of local laws and regulations"
$ws.Range("D16").Value = "5-14"
$ws.Range("E16").Value = "Missing"

$ws.Range("A17").Value = 85
$ws.Range("B17").Value = "Flight altitude exceeds nnn . Aircraft may be in violation of local laws and regulations . Check and make sure you have obtained proper authorization to fly in this airspace ."
$ws.Range("C17").Value = "Aircraft may be in violation of local laws"
$ws.Range("D17").Value = "5-12"
$ws.Range("E17").Value = "'False"

$ws.Range("A18").Value = 91
$ws.Range("B18").Value = "GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn ."
$ws.Range("C18").Value = "GEO Zone Info: The target area is in an Altitude Zone"
$ws.Range("D18").Value = "0-10"
$ws.Range("E18").Value = "Missing"

$ws.Range("A19").Value = 91
$ws.Range("B19").Value = "GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn ."
$ws.Range("C19").Value = "GEO Zone Info:"
$ws.Range("D19").Value = "0-2"
$ws.Range("E19").Value = "'False"

$ws.Range("A20").Value = 91
$ws.Range("B20").Value = "GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn ."
$ws.Range("C20").Value = "The target area is in an Altitude Zone"
$ws.Range("D20").Value = "3-10"
$ws.Range("E20").Value = "'False"

$ws.Range("A21").Value = 115
$ws.Range("B21").Value = "Landin ."
$ws.Range("C21").Value = "Landin"
$ws.Range("D21").Value = "0-0"
$ws.Range("E21").Value = "Missing"

$ws.Range("A22").Value = 141
$ws.Range("B22").Value = "Remote controller signal weak . Adjust remote controller antennas ."
$ws.Range("C22").Value = "Adjust remote controller antennas"
$ws.Range("D22").Value = "5-8"
$ws.Range("E22").Value = "Missing"

$ws.Range("A23").Value = 163
$ws.Range("B23").Value = "Warnin ."
$ws.Range("C23").Value = "Warnin"
$ws.Range("D23").Value = "0-0"
$ws.Range("E23").Value = "Missing"

$ws.Rows("24:33").Delete()
